$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.642.47"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.716.16"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +4.43%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.995"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.59%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "410.96"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.42"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.860.27"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +8.66%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.626"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.99%  "

$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.736"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.167"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000334"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +7.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.60"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.67%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.08"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.296.19"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.14%  "

$ws.Range("E16").Value = "  -1.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.779.78"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.87%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.13"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.02"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.10"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "65.821.41"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "424.91"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -6.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.74"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +12.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.11"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.02"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -6.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "36.52"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +6.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.21"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.57"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.25%  "

$ws.Range("E29").Value = "  +6.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.53"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.120"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.72"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.02"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.36%  "

$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.29"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.73%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.159"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.98"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0470"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -6.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.94"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +26.18%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.141"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.90%  "

$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0677"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -15.89%  "

$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.991"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.76%  "

$ws.Range("B43").Value = "LidoDAOToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.35"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.68%  "

$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "145.05"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.76%  "

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.07"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.60%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.30"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +20.35%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.16"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +23.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.26"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.81"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -7.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.52"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -9.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.292"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.85%  "
